$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 333.33334
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 360
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 360
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -586

$ws.Range("H132").Value = 19595.75
$ws.Range("I132").Value = 20173.371
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 60520.113
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -57990.113
$ws.Range("N132").Value = -17060

$ws.Range("H137").Value = 19232082
$ws.Range("I137").Value = 29412782
$ws.Range("J137").Value = 1871.1666
$ws.Range("K137").Value = 88238346
$ws.Range("L137").Value = 5613.4998
$ws.Range("M137").Value = -88235796
$ws.Range("N137").Value = -10713.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 2284.4
$ws.Range("I31").Value = 2284.4
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2284.4
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1990.4
$ws.Range("N31").ClearContents()

$ws.Range("H32").Value = 20405.014
$ws.Range("I32").Value = 4713.0156
$ws.Range("J32").Value = 145941
$ws.Range("K32").Value = 4713.0156
$ws.Range("L32").Value = 145941
$ws.Range("M32").Value = -4426.0156
$ws.Range("N32").Value = -146515

$ws.Range("H61").Value = 2921.5
$ws.Range("I61").Value = 2761.6155
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2761.6155
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2549.6155
$ws.Range("N61").Value = -5424

$ws.Range("H74").Value = 4197.5
$ws.Range("I74").Value = 1119.3103
$ws.Range("J74").Value = 12312.728
$ws.Range("K74").Value = 1119.3103
$ws.Range("L74").Value = 12312.728
$ws.Range("M74").Value = -245.3103000000001
$ws.Range("N74").Value = -14060.728

$ws.Range("H77").Value = 4197.5
$ws.Range("I77").Value = 1119.3103
$ws.Range("J77").Value = 12312.728
$ws.Range("K77").Value = 5596.5515
$ws.Range("L77").Value = 61563.64
$ws.Range("M77").Value = -1228.5515
$ws.Range("N77").Value = -70299.64

$ws.Range("H136").Value = 2921.5
$ws.Range("I136").Value = 2761.6155
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 8284.8465
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -5734.8465
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3234.3555
$ws.Range("I105").Value = 3032.3438
$ws.Range("J105").Value = 3731.6155
$ws.Range("K105").Value = 3032.3438
$ws.Range("L105").Value = 3731.6155
$ws.Range("M105").Value = -1285.3438
$ws.Range("N105").Value = -7225.6155

$ws.Range("H106").Value = 30000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 30000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524

$ws.Range("H132").Value = 46086.668
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 46086.668
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 46086.668
$ws.Range("N132").Value = -56206.668

$ws.Range("H134").Value = 26319006
$ws.Range("I134").Value = 38463756
$ws.Range("J134").Value = 5384.1665
$ws.Range("K134").Value = 115391268
$ws.Range("L134").Value = 16152.4995
$ws.Range("M134").Value = -115388733
$ws.Range("N134").Value = -21222.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1663.5807
$ws.Range("I58").Value = 1136.6
$ws.Range("J58").Value = 3859.3333
$ws.Range("K58").Value = 1136.6
$ws.Range("L58").Value = 3859.3333
$ws.Range("M58").Value = -933.5999999999999
$ws.Range("N58").Value = -4265.3333

$ws.Range("H105").Value = 1055.5555
$ws.Range("I105").Value = 1087.1428
$ws.Range("J105").Value = 945
$ws.Range("K105").Value = 1087.1428
$ws.Range("L105").Value = 945
$ws.Range("M105").Value = 659.8571999999999
$ws.Range("N105").Value = -4439

$ws.Range("H134").Value = 3594.7058
$ws.Range("I134").Value = 2240.3
$ws.Range("J134").Value = 5529.5713
$ws.Range("K134").Value = 6720.900000000001
$ws.Range("L134").Value = 16588.7139
$ws.Range("M134").Value = -4185.900000000001
$ws.Range("N134").Value = -21658.7139

$ws.Range("H136").Value = 1663.5807
$ws.Range("I136").Value = 1136.6
$ws.Range("J136").Value = 3859.3333
$ws.Range("K136").Value = 3409.8
$ws.Range("L136").Value = 11577.9999
$ws.Range("M136").Value = -859.7999999999997
$ws.Range("N136").Value = -16677.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 3442.6
$ws.Range("I82").Value = 1013
$ws.Range("J82").Value = 4050
$ws.Range("K82").Value = 3039
$ws.Range("L82").Value = 12150
$ws.Range("M82").Value = -2633
$ws.Range("N82").Value = -12962

$ws.Range("H85").Value = 3442.6
$ws.Range("I85").Value = 1013
$ws.Range("J85").Value = 4050
$ws.Range("K85").Value = 3039
$ws.Range("L85").Value = 12150
$ws.Range("M85").Value = -1635
$ws.Range("N85").Value = -14958

$ws.Range("H131").Value = 9805773
$ws.Range("I131").Value = 575
$ws.Range("J131").Value = 11113133
$ws.Range("K131").Value = 1725
$ws.Range("L131").Value = 33339399
$ws.Range("M131").Value = 3315
$ws.Range("N131").Value = -33349479

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 650
$ws.Range("I22").Value = 375
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 375
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -80
$ws.Range("N22").Value = -1790

$ws.Range("H27").Value = 650
$ws.Range("I27").Value = 375
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 375
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = -268
$ws.Range("N27").Value = -1414

$ws.Range("H132").Value = 2698.1667
$ws.Range("I132").Value = 1798.0222
$ws.Range("J132").Value = 5398.6
$ws.Range("K132").Value = 5394.0666
$ws.Range("L132").Value = 16195.8
$ws.Range("M132").Value = -2864.0666
$ws.Range("N132").Value = -21255.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 30000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 30000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

$ws.Range("H95").Value = 30000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 30000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 30000
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -35492

$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H125").Value = 29369.375
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 29369.375
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 29369.375
$ws.Range("N125").Value = -39209.375

$ws.Range("H128").Value = 42857.5
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 42857.5
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 42857.5
$ws.Range("N128").Value = -52817.5

$ws.Range("H131").Value = 74750
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 74750
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 74750
$ws.Range("N131").Value = -84830

